$d = $word.ActiveDocument

$replacements = @(
    @("29×14=", "75×77="),
    @("38×44=", "68×35="),
    @("21×46=", "69×26="),
    @("96×85=", "61×53="),
    @("75×65=", "32×57="),
    @("87×51=", "90×46="),
    @("68×15=", "34×62="),
    @("38×80=", "45×18="),
    @("43×63=", "95×40="),
    @("49×19=", "78×26="),
    @("51×71=", "72×56="),
    @("67×63=", "51×18="),
    @("81×77=", "54×75="),
    @("20×52=", "38×49="),
    @("38×13=", "75×19="),
    @("22×89=", "95×16="),
    @("94×30=", "83×66="),
    @("74×62=", "79×79="),
    @("27×54=", "38×45="),
    @("26×70=", "97×48="),
    @("36×77=", "68×21="),
    @("84×32=", "88×93="),
    @("87×73=", "94×11="),
    @("61×36=", "88×17="),
    @("29×70=", "38×77=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
